# Add three new data rows (65, 66, 67) to Sheet1, continuing the daily
# series with the same B:J values as the preceding rows and consecutive
# date serials in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 65; Date = 45621 },
    @{ Row = 66; Date = 45622 },
    @{ Row = 67; Date = 45623 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Range("A$r").Value = $entry.Date
    $ws.Range("B$r").Value = 116.4121952
    $ws.Range("C$r").Value = 0.00170247
    $ws.Range("D$r").Value = 0.008850780000000001
    $ws.Range("E$r").Value = 0.06933635
    $ws.Range("F$r").Value = 12792.90181321
    $ws.Range("G$r").Value = 465.80531254
    $ws.Range("H$r").Value = 0.24
    $ws.Range("I$r").Value = 1.7904431
    $ws.Range("J$r").Value = 485.38834923

    # Mirror the date cell's formatting (bold, bordered, centered,
    # YYYY-MM-DD HH:MM:SS) from the row directly above it, reusing the
    # existing style rather than creating a new one.
    $ws.Range("A64").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
}
